$d = $word.ActiveDocument

# 1. Update the supervisor's name.
#    This run sits right next to another run ("Преподаватель: ") that shares
#    identical formatting. A plain Find/Replace would cause the runtime to
#    coalesce the two runs into one, but the target keeps them as two
#    separate runs. Toggling a throwaway formatting property (and restoring
#    it) keeps the edited run distinct so it is not merged with its sibling.
$nameRange = $d.Content
$found = $nameRange.Find.Execute("Черемисинов Максим", $true, $false, $false, $false, $false,
                                  $true, 1, $false, "", 0)
if ($found) {
    $nameRange.Font.Bold = 1
    $nameRange.Text = "Миронов Евгений Сергеевич"
    $nameRange.Font.Bold = 0
}

# 2. Coalesce the run boundaries around "Count of threads: N  Программа работала ..."
#    lines so that the split runs become single merged runs (same visible text).
$d.Content.Find.Execute(": 10 ", $true, $false, $false, $false, $false,
                         $true, 1, $false, ": 10 ", 2)

$d.Content.Find.Execute(": 20 ", $true, $false, $false, $false, $false,
                         $true, 1, $false, ": 20 ", 2)

$d.Content.Find.Execute(": 50 ", $true, $false, $false, $false, $false,
                         $true, 1, $false, ": 50 ", 2)

$d.Content.Find.Execute(": 100 ", $true, $false, $false, $false, $false,
                         $true, 1, $false, ": 100 ", 2)
